# Rename risk driver 2 (cell A3) to its new label.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Suitability for circular proposition"

# Mirror the natural cursor movement a user gets after typing into A3
# and pressing Enter (selection lands on the cell below).
$ws.Range("A4").Select()
